$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "Devem existir máquinas alocadas no cl" + bookmark(_GoBack) +
#           "iente em questão [Caso de uso 51]"
#  -> single run "Devem existir máquinas alocadas no cliente em questão [Caso de uso 51]"
# A Find/Replace spanning the whole phrase merges the two runs (identical
# rPr) into one and drops the bookmark that sat between them.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Devem existir máquinas alocadas no cliente em questão [Caso de uso 51]",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Devem existir máquinas alocadas no cliente em questão [Caso de uso 51]", 2)

# ---------------------------------------------------------------------
# Change 2: " [9.1]" (single run) -> " [10" + ".1]" (two runs, same rPr)
# This run shares rPr with its left neighbour ("”."), so any text edit
# would otherwise merge them. Isolate the run (temporarily flip Bold on)
# before editing the text, then split the run in two by toggling Bold
# on/off (format-only op => forces a run split without re-merging
# neighbours) right on the boundary between "10" and ".1]", and finally
# drop Bold back off both halves.
# ---------------------------------------------------------------------
$r2pre = $d.Content
$r2pre.Find.Execute(" [9.1]")
$full2 = $d.Range($r2pre.Start, $r2pre.End)
$full2.Bold = $true

$d.Content.Find.Execute(" [9.1]", $true, $false, $false, $false, $false, $true, 1, $false, " [10.1]", 2)

$r2 = $d.Content
$r2.Find.Execute("[10.1]")
$splitPoint = $r2.Start + 3
$secondHalf = $d.Range($splitPoint, $r2.End)
$secondHalf.Bold = $true
$secondHalf.Bold = $false

$firstHalf = $d.Range($r2.Start - 1, $splitPoint)
$firstHalf.Bold = $false

# ---------------------------------------------------------------------
# Change 3: "Falha ao acessar servidor" -> "Erro ao acessar banco de dados"
# This run shares rPr with its left neighbours ("2" and ".1. "), so we
# isolate it (temporarily drop Bold) before editing the text and restore
# Bold afterwards (format-only ops don't trigger the run-merge pass).
# ---------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("Falha ao acessar servidor")
$full3 = $d.Range($r3.Start, $r3.End)
$full3.Bold = $false
$d.Content.Find.Execute("Falha ao acessar servidor", $true, $false, $false, $false, $false, $true, 1, $false, "Erro ao acessar banco de dados", 2)
$r3b = $d.Content
$r3b.Find.Execute("Erro ao acessar banco de dados")
$full3b = $d.Range($r3b.Start, $r3b.End)
$full3b.Bold = $true

# ---------------------------------------------------------------------
# Change 4: "9" -> "10"; "O" -> "Ator" (+ new _GoBack bookmark right after
# it); "perador não preenche o motivo" -> " não preenche o motivo"
# All four runs in this paragraph ("9" / ".1. " / "O" / "perador...")
# share identical rPr. We isolate the ones we must edit (run1 "9" and
# run3 "O") by flipping Bold off, edit all three texts while isolated,
# then restore Bold on runs 1 and 3 with format-only ops so nothing
# re-merges, and finally drop a new _GoBack bookmark between "Ator" and
# " não preenche o motivo".
# ---------------------------------------------------------------------
$anchor4 = $d.Content
$anchor4.Find.Execute("9.1. Operador não preenche o motivo")
$s = $anchor4.Start

$run1 = $d.Range($s, $s + 1)
$run1.Bold = $false
$run3 = $d.Range($s + 5, $s + 6)
$run3.Bold = $false

$run1edit = $d.Range($s, $s + 1)
$run1edit.Text = "10"

$delta = 1
$run3edit = $d.Range($s + 5 + $delta, $s + 6 + $delta)
$run3edit.Text = "Ator"

$delta2 = 3
$run4edit = $d.Range($s + 6 + $delta + $delta2, $s + 36 + $delta + $delta2)
$run4edit.Text = " não preenche o motivo"

$run1restore = $d.Range($s, $s + 2)
$run1restore.Bold = $true
$run3restore = $d.Range($s + 6, $s + 10)
$run3restore.Bold = $true

$bmPos = $d.Range($s + 10, $s + 10)
$d.Bookmarks.Add("_GoBack", $bmPos)

# ---------------------------------------------------------------------
# Change 5: "Fim caso de uso" + "." (two runs, third occurrence of the
# phrase in the document) -> single run "Sistema retorna para o passo 9."
# We scope the search to start right after the preceding sentence so we
# hit the correct occurrence (the phrase repeats earlier in the doc).
# ---------------------------------------------------------------------
$anchor5 = $d.Content
$anchor5.Find.Execute("Sistema exibe mensagem informando que o campo deve ser preenchido.")
$scope5 = $d.Range($anchor5.End, $d.Content.End)
$scope5.Find.Execute("Fim caso de uso.")
$target5 = $d.Range($scope5.Start, $scope5.End)
$target5.Text = "Sistema retorna para o passo 9."
